# The sheet has a block of "label" columns (Alain / Henri / Tony / Dulcinee
# repeated, then OUI/NON answer columns) that used to stop right before the
# trailing "Adresse de courriel" / empty columns (ANU:ANV). The commit
# removes the label-boolean-by-name approach and extends the columns grid by
# one more repeat-group set (16 columns), so the data must now be addressed
# by column letters. Concretely: insert 16 new columns at ANU (shifting the
# trailing email / empty columns from ANU:ANV out to AOK:AOL), and populate
# the newly inserted columns with a copy of the previous 16-column block
# (ANE:ANT), exactly continuing the repeating pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$shiftRight = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight
$pasteAll = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteAll

# Insert 16 blank columns at ANU:AOJ, pushing the old ANU:ANV (email +
# trailing empty column) out to AOK:AOL.
$ws.Range("ANU1:AOJ1").EntireColumn.Insert($shiftRight)

# Fill the newly inserted columns with a copy of the previous 16-column
# repeating block (ANE:ANT), covering every used row (header row 1 plus the
# 8 data rows).
$ws.Range("ANE1:ANT9").Copy()
$ws.Range("ANU1").PasteSpecial($pasteAll)
